$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the two new algorithm rows with their computed metrics
$ws.Range("A4").Value = "Support Vector Regressor"
$ws.Range("B4").Value = 0.2636
$ws.Range("C4").Value = 129.03
$ws.Range("D4").Value = 31605.17
$ws.Range("E4").Value = 177.78

$ws.Range("A5").Value = "K Nearest Neighbors"
$ws.Range("B5").Value = 0.8458
$ws.Range("C5").Value = 57.6242
$ws.Range("D5").Value = 6619.582
$ws.Range("E5").Value = 81.3608

$ws.Range("D4").Select()
